$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Back-Biceps-Abs": insert a new "Back Hyperextension" row before the
# existing "Planks" row (row 9), shifting Planks/Leg Raises/Sit Ups down.
# ---------------------------------------------------------------------------
$wsBack = $wb.Worksheets.Item("Back-Biceps-Abs")
[void]$wsBack.Rows.Item(9).Insert()
$wsBack.Range("A9").Value = "Back Hyperextension"
$wsBack.Range("B9").Value = 4
$wsBack.Range("C9").Value = 12

# ---------------------------------------------------------------------------
# Sheet "Chest-Triceps": row 7 used to duplicate "Bench Press" - rename it to
# "Dumbbell Press". Row 8 ("Ab Roller") rep count changes from 6 to 10.
# ---------------------------------------------------------------------------
$wsChest = $wb.Worksheets.Item("Chest-Triceps")
$wsChest.Range("A7").Value = "Dumbbell Press"
$wsChest.Range("C8").Value = 10

# ---------------------------------------------------------------------------
# Sheet "Legs": no exercise/value changes, only the remembered selection
# moves from A5 to A6 (handled in the selection section below).
# ---------------------------------------------------------------------------
$wsLegs = $wb.Worksheets.Item("Legs")

# ---------------------------------------------------------------------------
# Sheet "Shoulders-Traps": remove the "Dumbbell Shrugs" row entirely; the
# "Farmers Carry" row that followed shifts up to take its place.
# ---------------------------------------------------------------------------
$wsShoulders = $wb.Worksheets.Item("Shoulders-Traps")
[void]$wsShoulders.Rows.Item(9).Delete()

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping to mirror the saved UI state.
# ---------------------------------------------------------------------------
[void]$wsBack.Range("A1").Select()
[void]$wsLegs.Range("A6").Select()
[void]$wsShoulders.Range("A1").Select()

[void]$wsChest.Activate()
[void]$wsChest.Range("C8").Select()
